# The "Punny signs" slide (7th slide in the deck, sldId=257) is removed
# entirely. The slide that follows it ("THANK YOU !!!", sldId=263) shifts
# up to become the new (final) 7th slide.
$p = $ppt.ActivePresentation
$p.Slides.Item(7).Delete()
